$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.04599157227274
$ws.Cells.Item(2, 4).Value = 1.047912524197416
$ws.Cells.Item(2, 5).Value = 1.043620486007928
$ws.Cells.Item(2, 6).Value = 1.053562309972952
$ws.Cells.Item(2, 9).Value = 1.043376629930183
$ws.Cells.Item(2, 10).Value = 1.05104835039516
$ws.Cells.Item(2, 11).Value = 1.050673836319028
$ws.Cells.Item(2, 12).Value = 1.046393837905036
$ws.Cells.Item(2, 13).Value = 1.056307935735965

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.047236337949261
$ws.Cells.Item(3, 4).Value = 1.048887649499187
$ws.Cells.Item(3, 5).Value = 1.044690870434723
$ws.Cells.Item(3, 6).Value = 1.055000161718333
$ws.Cells.Item(3, 9).Value = 1.043777781630567
$ws.Cells.Item(3, 10).Value = 1.051939546910758
$ws.Cells.Item(3, 11).Value = 1.051460478548694
$ws.Cells.Item(3, 12).Value = 1.04727461669178
$ws.Cells.Item(3, 13).Value = 1.057557257706075

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.048041143505223
$ws.Cells.Item(4, 4).Value = 1.049517997423588
$ws.Cells.Item(4, 5).Value = 1.045383153643662
$ws.Cells.Item(4, 6).Value = 1.055930257446881
$ws.Cells.Item(4, 9).Value = 1.044035803055177
$ws.Cells.Item(4, 10).Value = 1.052515067234555
$ws.Cells.Item(4, 11).Value = 1.051968250394145
$ws.Cells.Item(4, 12).Value = 1.047843618625594
$ws.Cells.Item(4, 13).Value = 1.058364843853478

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.048379333379124
$ws.Cells.Item(5, 4).Value = 1.049782848163068
$ws.Cells.Item(5, 5).Value = 1.045674113111397
$ws.Cells.Item(5, 6).Value = 1.056321204275614
$ws.Cells.Item(5, 9).Value = 1.044143905090804
$ws.Cells.Item(5, 10).Value = 1.05275674427203
$ws.Cells.Item(5, 11).Value = 1.052181422931124
$ws.Cells.Item(5, 12).Value = 1.048082608404508
$ws.Cells.Item(5, 13).Value = 1.058704163072175

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.048436108218339
$ws.Cells.Item(6, 4).Value = 1.049827309165889
$ws.Cells.Item(6, 5).Value = 1.045722962062032
$ws.Cells.Item(6, 6).Value = 1.056386842252297
$ws.Cells.Item(6, 9).Value = 1.044162034221213
$ws.Cells.Item(6, 10).Value = 1.052797307035894
$ws.Cells.Item(6, 11).Value = 1.052217198303267
$ws.Cells.Item(6, 12).Value = 1.048122723071852
$ws.Cells.Item(6, 13).Value = 1.05876112522525

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.048045663001911
$ws.Cells.Item(7, 4).Value = 1.049521536949698
$ws.Cells.Item(7, 5).Value = 1.045387041756104
$ws.Cells.Item(7, 6).Value = 1.055935481549429
$ws.Cells.Item(7, 9).Value = 1.044037248973103
$ws.Cells.Item(7, 10).Value = 1.052518297600056
$ws.Cells.Item(7, 11).Value = 1.05197109997206
$ws.Cells.Item(7, 12).Value = 1.047846812875543
$ws.Cells.Item(7, 13).Value = 1.058369378595552

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.046412381262706
$ws.Cells.Item(8, 4).Value = 1.048242202455471
$ws.Cells.Item(8, 5).Value = 1.043982296316305
$ws.Cells.Item(8, 6).Value = 1.054048301340043
$ws.Cells.Item(8, 9).Value = 1.043512522700286
$ws.Cells.Item(8, 10).Value = 1.051349772405231
$ws.Cells.Item(8, 11).Value = 1.050939943162425
$ws.Cells.Item(8, 12).Value = 1.046691692980654
$ws.Cells.Item(8, 13).Value = 1.056730319567519

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.043529281480873
$ws.Cells.Item(9, 4).Value = 1.04598300060436
$ws.Cells.Item(9, 5).Value = 1.041504364036813
$ws.Cells.Item(9, 6).Value = 1.050720438655875
$ws.Cells.Item(9, 9).Value = 1.042575966045973
$ws.Cells.Item(9, 10).Value = 1.049281833199636
$ws.Cells.Item(9, 11).Value = 1.049113352797281
$ws.Cells.Item(9, 12).Value = 1.044649089806824
$ws.Cells.Item(9, 13).Value = 1.053835732552515

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041603624290201
$ws.Cells.Item(10, 4).Value = 1.044473494456281
$ws.Cells.Item(10, 5).Value = 1.039850543879658
$ws.Cells.Item(10, 6).Value = 1.048499982498131
$ws.Cells.Item(10, 9).Value = 1.041943512934833
$ws.Cells.Item(10, 10).Value = 1.047897128045559
$ws.Cells.Item(10, 11).Value = 1.047889088927991
$ws.Cells.Item(10, 12).Value = 1.043282441092064
$ws.Cells.Item(10, 13).Value = 1.051901506647815

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040768897487659
$ws.Cells.Item(11, 4).Value = 1.043819036668174
$ws.Cells.Item(11, 5).Value = 1.039133950918332
$ws.Cells.Item(11, 6).Value = 1.047537995765653
$ws.Cells.Item(11, 9).Value = 1.041667721625165
$ws.Cells.Item(11, 10).Value = 1.04729606403219
$ws.Cells.Item(11, 11).Value = 1.047357393743446
$ws.Cells.Item(11, 12).Value = 1.042689475146048
$ws.Cells.Item(11, 13).Value = 1.051062846536271

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040458702981133
$ws.Cells.Item(12, 4).Value = 1.043575814683647
$ws.Cells.Item(12, 5).Value = 1.038867702064533
$ws.Cells.Item(12, 6).Value = 1.047180588629667
$ws.Cells.Item(12, 9).Value = 1.041564988375986
$ws.Cells.Item(12, 10).Value = 1.047072577475409
$ws.Cells.Item(12, 11).Value = 1.047159658735372
$ws.Cells.Item(12, 12).Value = 1.04246903881932
$ws.Cells.Item(12, 13).Value = 1.050751156925892

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040525247145238
$ws.Cells.Item(13, 4).Value = 1.043627992418707
$ws.Cells.Item(13, 5).Value = 1.0389248167513
$ws.Cells.Item(13, 6).Value = 1.047257257486262
$ws.Cells.Item(13, 9).Value = 1.041587038245915
$ws.Cells.Item(13, 10).Value = 1.047120526307462
$ws.Cells.Item(13, 11).Value = 1.047202084450761
$ws.Cells.Item(13, 12).Value = 1.042516331450811
$ws.Cells.Item(13, 13).Value = 1.050818023347412

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040743259587657
$ws.Cells.Item(14, 4).Value = 1.043798934471559
$ws.Cells.Item(14, 5).Value = 1.039111944235785
$ws.Cells.Item(14, 6).Value = 1.047508454091579
$ws.Cells.Item(14, 9).Value = 1.041659235628998
$ws.Cells.Item(14, 10).Value = 1.047277595161157
$ws.Cells.Item(14, 11).Value = 1.047341053809663
$ws.Cells.Item(14, 12).Value = 1.042671257544681
$ws.Cells.Item(14, 13).Value = 1.051037085749549

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040877565605864
$ws.Cells.Item(15, 4).Value = 1.043904240585668
$ws.Cells.Item(15, 5).Value = 1.039227229731005
$ws.Cells.Item(15, 6).Value = 1.047663213534099
$ws.Cells.Item(15, 9).Value = 1.041703680075249
$ws.Cells.Item(15, 10).Value = 1.047374340622472
$ws.Cells.Item(15, 11).Value = 1.047426645582852
$ws.Cells.Item(15, 12).Value = 1.042766688386729
$ws.Cells.Item(15, 13).Value = 1.051172034146484

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041659002968162
$ws.Cells.Item(16, 4).Value = 1.044516910948269
$ws.Cells.Item(16, 5).Value = 1.039898091510503
$ws.Cells.Item(16, 6).Value = 1.048563815111221
$ws.Cells.Item(16, 9).Value = 1.041961775411645
$ws.Cells.Item(16, 10).Value = 1.047936987329757
$ws.Cells.Item(16, 11).Value = 1.04792434228885
$ws.Cells.Item(16, 12).Value = 1.043321768844831
$ws.Cells.Item(16, 13).Value = 1.051957141617491

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.042148933055577
$ws.Cells.Item(17, 4).Value = 1.044900998618379
$ws.Cells.Item(17, 5).Value = 1.040318775653325
$ws.Cells.Item(17, 6).Value = 1.049128597312489
$ws.Cells.Item(17, 9).Value = 1.04212315279065
$ws.Cells.Item(17, 10).Value = 1.048289523306851
$ws.Cells.Item(17, 11).Value = 1.048236109468455
$ws.Cells.Item(17, 12).Value = 1.043669633725673
$ws.Cells.Item(17, 13).Value = 1.052449313744735

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042434613978413
$ws.Cells.Item(18, 4).Value = 1.045124950369871
$ws.Cells.Item(18, 5).Value = 1.040564107656083
$ws.Cells.Item(18, 6).Value = 1.049457975915348
$ws.Cells.Item(18, 9).Value = 1.042217094885659
$ws.Cells.Item(18, 10).Value = 1.04849500919587
$ws.Cells.Item(18, 11).Value = 1.048417805613137
$ws.Cells.Item(18, 12).Value = 1.043872421994763
$ws.Cells.Item(18, 13).Value = 1.052736280935434

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042532009149967
$ws.Cells.Item(19, 4).Value = 1.045201298613522
$ws.Cells.Item(19, 5).Value = 1.040647751812387
$ws.Cells.Item(19, 6).Value = 1.049570277299197
$ws.Cells.Item(19, 9).Value = 1.042249095103076
$ws.Cells.Item(19, 10).Value = 1.048565050502515
$ws.Cells.Item(19, 11).Value = 1.04847973352451
$ws.Cells.Item(19, 12).Value = 1.043941548046484
$ws.Cells.Item(19, 13).Value = 1.052834111066372

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.042096377223142
$ws.Cells.Item(20, 4).Value = 1.0448597979477
$ws.Cells.Item(20, 5).Value = 1.040273644997239
$ws.Cells.Item(20, 6).Value = 1.049068006665269
$ws.Cells.Item(20, 9).Value = 1.042105857836546
$ws.Cells.Item(20, 10).Value = 1.048251714280189
$ws.Cells.Item(20, 11).Value = 1.048202675578729
$ws.Cells.Item(20, 12).Value = 1.043632323073829
$ws.Cells.Item(20, 13).Value = 1.052396519588384

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040679064247993
$ws.Cells.Item(21, 4).Value = 1.043748599827903
$ws.Cells.Item(21, 5).Value = 1.039056841927581
$ws.Cells.Item(21, 6).Value = 1.047434485268979
$ws.Cells.Item(21, 9).Value = 1.041637983372981
$ws.Cells.Item(21, 10).Value = 1.047231348520421
$ws.Cells.Item(21, 11).Value = 1.047300137433293
$ws.Cells.Item(21, 12).Value = 1.042625640725031
$ws.Cells.Item(21, 13).Value = 1.050972582168968

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039787132505778
$ws.Cells.Item(22, 4).Value = 1.043049208261572
$ws.Cells.Item(22, 5).Value = 1.038291358282753
$ws.Cells.Item(22, 6).Value = 1.04640694605967
$ws.Cells.Item(22, 9).Value = 1.041342121574375
$ws.Cells.Item(22, 10).Value = 1.046588503173409
$ws.Cells.Item(22, 11).Value = 1.04673128768493
$ws.Cells.Item(22, 12).Value = 1.041991643237431
$ws.Cells.Item(22, 13).Value = 1.050076287800698

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040260040563807
$ws.Cells.Item(23, 4).Value = 1.043420039604578
$ws.Cells.Item(23, 5).Value = 1.038697197369052
$ws.Cells.Item(23, 6).Value = 1.046951711200015
$ws.Cells.Item(23, 9).Value = 1.041499124216818
$ws.Cells.Item(23, 10).Value = 1.046929411847509
$ws.Cells.Item(23, 11).Value = 1.047032978034069
$ws.Cells.Item(23, 12).Value = 1.042327838227849
$ws.Cells.Item(23, 13).Value = 1.050551527496223

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.042120125210276
$ws.Cells.Item(24, 4).Value = 1.044878415004761
$ws.Cells.Item(24, 5).Value = 1.04029403774002
$ws.Cells.Item(24, 6).Value = 1.049095385122726
$ws.Cells.Item(24, 9).Value = 1.042113673258382
$ws.Cells.Item(24, 10).Value = 1.048268798991731
$ws.Cells.Item(24, 11).Value = 1.048217783385121
$ws.Cells.Item(24, 12).Value = 1.043649182507991
$ws.Cells.Item(24, 13).Value = 1.052420375330119

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.044275250331213
$ws.Cells.Item(25, 4).Value = 1.046567644590186
$ws.Cells.Item(25, 5).Value = 1.042145288636674
$ws.Cells.Item(25, 6).Value = 1.051581083309771
$ws.Cells.Item(25, 9).Value = 1.04281950801395
$ws.Cells.Item(25, 10).Value = 1.04981750696802
$ws.Cells.Item(25, 11).Value = 1.049586714336185
$ws.Cells.Item(25, 12).Value = 1.045178009103105
$ws.Cells.Item(25, 13).Value = 1.054584829280424
